# Deep sea double count fix
# Updates landings / weighted-% figures for "Sharks" (and derived "Global")
# rows on the "Status by Landings (Area)" and "Status by Landings (Tier)"
# sheets to correct a double-counting bug.

$wb = $excel.ActiveWorkbook

# --- Sheet: "Status by Landings (Area)" -----------------------------------
$wsArea = $wb.Worksheets.Item("Status by Landings (Area)")

$wsArea.Range("C2").Value = 0.02678414
$wsArea.Range("C3").Value = 0.02084913
$wsArea.Range("C4").Value = 0.00601671
$wsArea.Range("C5").Value = 0.04763327000000001
$wsArea.Range("C6").Value = 0.00601671
$wsArea.Range("C7").Value = 49.92385831271513
$wsArea.Range("C8").Value = 38.86139379735091
$wsArea.Range("C9").Value = 11.21474788993398
$wsArea.Range("C10").Value = 88.78525211006604
$wsArea.Range("C11").Value = 11.21474788993398

# --- Sheet: "Status by Landings (Tier)" ------------------------------------
$wsTier = $wb.Worksheets.Item("Status by Landings (Tier)")

# Row 4 - "Sharks"
$wsTier.Range("B4").Value = 0.02678414
$wsTier.Range("C4").Value = 0.02084913
$wsTier.Range("D4").Value = 0.00601671
$wsTier.Range("E4").Value = 0.04763327000000001
$wsTier.Range("F4").Value = 0.00601671
$wsTier.Range("G4").Value = 49.92385831271513
$wsTier.Range("H4").Value = 38.86139379735091
$wsTier.Range("I4").Value = 11.21474788993398
$wsTier.Range("J4").Value = 88.78525211006604
$wsTier.Range("K4").Value = 11.21474788993398

# Row 5 - "Global"
$wsTier.Range("B5").Value = 0.02678414
$wsTier.Range("C5").Value = 0.02084913
$wsTier.Range("D5").Value = 0.00601671
$wsTier.Range("E5").Value = 0.04763327000000001
$wsTier.Range("F5").Value = 0.00601671
$wsTier.Range("G5").Value = 49.92385831271513
$wsTier.Range("H5").Value = 38.86139379735091
$wsTier.Range("I5").Value = 11.21474788993398
$wsTier.Range("J5").Value = 88.78525211006604
$wsTier.Range("K5").Value = 11.21474788993398
